# Apply the "added new redline channels" edit to the channels sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("channels")

# --- Fix up the absolute-reference formulas in rows 11-14 ---
$ws.Range("C11").Formula = '=$C$3'
$ws.Range("D11").Formula = '=$D$3'
$ws.Range("C12").Formula = '=$C$3'
$ws.Range("D12").Formula = '=$D$3'
$ws.Range("C13").Formula = '=$C$2'
$ws.Range("D13").Formula = '=$D$2'
$ws.Range("C14").Formula = '=$C$2'
$ws.Range("D14").Formula = '=$D$2'

# --- New rows 15-18: redline channels (inserted ahead of the old 15/16 rows,
# which shift down to 19/20). Written in document order so the shared
# string table gets the new strings allocated right after OX_LOWER_SETP. ---

# --- Row 15: FU_UPPER_REDLINE ---
$ws.Range("A15").Value = "FU_UPPER_REDLINE"
$ws.Range("B15").Value = 18
$ws.Range("C15").Formula = '=$C$3'
$ws.Range("D15").Formula = '=$D$3'
$ws.Range("E15").Value = "f64"

# --- Row 16: FU_LOWER_REDLINE ---
$ws.Range("A16").Value = "FU_LOWER_REDLINE"
$ws.Range("B16").Value = 19
$ws.Range("C16").Formula = '=$C$3'
$ws.Range("D16").Formula = '=$D$3'
$ws.Range("E16").Value = "f64"

# --- Row 17: OX_UPPER_REDLINE ---
$ws.Range("A17").Value = "OX_UPPER_REDLINE"
$ws.Range("B17").Value = 20
$ws.Range("C17").Formula = '=$C$2'
$ws.Range("D17").Formula = '=$D$2'
$ws.Range("E17").Value = "f64"

# --- Row 18: OX_LOWER_REDLINE ---
$ws.Range("A18").Value = "OX_LOWER_REDLINE"
$ws.Range("B18").Value = 21
$ws.Range("C18").Formula = '=$C$2'
$ws.Range("D18").Formula = '=$D$2'
$ws.Range("E18").Value = "f64"

# --- Preserve / relocate the two rows that used to sit at 15-16 ---
# They move down to rows 19-20 to make room for the 4 new REDLINE rows.
$ws.Range("A19").Value = "FREE_SPACE"
$ws.Range("B19").Value = 22
$ws.Range("E19").Value = "u64"

$ws.Range("A20").Value = "random_ass_adc_channel"
$ws.Range("B20").Value = 6
$ws.Range("E20").Value = "u64"

# --- Selection moves to B20 to match the author's final cursor position ---
$ws.Range("B20").Select()
